# daily auto push: 2026-01-11 13:39 UTC
# Insert one new daily-ranking row (2026/01/11, 日, 19, 142) above the
# existing row 629, shifting the 2026/12/29 .. 2027/01/05 block down by
# one row (rows 629-670 -> 630-671) and extending the sheet dimension to
# A1:D671.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 629

# Push everything from row 629 down by one row.
$ws.Rows.Item($newRow).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/11"), not
# real Excel date serials, so force text interpretation before writing
# the value - then drop back to the default "Normal" style so the new
# cell matches its unstyled siblings.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/01/11"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "日"
$ws.Cells.Item($newRow, 3).Value = 19
$ws.Cells.Item($newRow, 4).Value = 142
